$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-ambiguous cell updates (Coin names, URLs, Volume percentages, and
# Price strings that are not valid numeric literals e.g. "61.508.96")
$ws.Range("D2").Value = "62.364.27"
$ws.Range("E2").Value = "  +12.83%  "
$ws.Range("D3").Value = "2.700.87"
$ws.Range("E3").Value = "  +14.77%  "
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("E5").Value = "  +9.14%  "
$ws.Range("E6").Value = "  +10.93%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "2.703.26"
$ws.Range("E9").Value = "  +14.62%  "
$ws.Range("E10").Value = "  +11.24%  "
$ws.Range("E11").Value = "  +13.05%  "
$ws.Range("E12").Value = "  +8.14%  "
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "3.103.76"
$ws.Range("E14").Value = "  +12.33%  "
$ws.Range("D15").Value = "61.690.75"
$ws.Range("E15").Value = "  +11.71%  "
$ws.Range("E16").Value = "  +13.70%  "
$ws.Range("E17").Value = "  +10.88%  "
$ws.Range("D18").Value = "2.674.06"
$ws.Range("E18").Value = "  +13.55%  "
$ws.Range("E19").Value = "  +6.47%  "
$ws.Range("E20").Value = "  +16.74%  "
$ws.Range("E21").Value = "  +11.81%  "
$ws.Range("E22").Value = "  +10.48%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +8.31%  "
$ws.Range("E25").Value = "  +9.18%  "
$ws.Range("E26").Value = "  +11.96%  "
$ws.Range("D27").Value = "2.738.45"
$ws.Range("E27").Value = "  +11.76%  "
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "0.0₃0877"
$ws.Range("E29").Value = "  +18.19%  "
$ws.Range("E30").Value = "  +8.65%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  +9.91%  "
$ws.Range("E33").Value = "  +9.30%  "
$ws.Range("E34").Value = "  +8.90%  "
$ws.Range("E35").Value = "  +11.41%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E36").Value = "  +12.92%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E37").Value = "  +12.58%  "
$ws.Range("E38").Value = "  +8.68%  "
$ws.Range("E39").Value = "  +15.28%  "
$ws.Range("E40").Value = "  +24.45%  "
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("E41").Value = "  +37.62%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E42").Value = "  +12.96%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E43").Value = "  +6.70%  "
$ws.Range("E44").Value = "  +13.15%  "
$ws.Range("E45").Value = "  +13.98%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +22.41%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E49").Value = "  +15.07%  "
$ws.Range("E50").Value = "  +9.10%  "
$ws.Range("D51").Value = "2.052.80"
$ws.Range("E51").Value = "  +14.37%  "

# Price cells whose text looks like a valid number (e.g. "0.995", "1.00").
# These must be forced to Text format first, otherwise Excel auto-converts them
# to numeric values and strips formatting such as trailing zeros.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.995"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.617"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.993"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "366.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.172"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.983"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "309.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.853"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.654"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0589"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.103"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.08"
$ws.Range("D49").Style = "Normal"
